{"js": "// Office.js (Word JavaScript API) script\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---- Helper: replace the text of an existing paragraph, preserving its\n// existing run formatting (font/size/etc. carry over automatically when\n// using InsertLocation.replace on the paragraph's own range). ----\nfunction setParaText(p, text) {\n  p.insertText(text, Word.InsertLocation.replace);\n}\n\nconst p = paragraphs.items;\n\n// 1) \"- Settings General save updates store via backend API.\"\n//    -> \"- Landing page Emergent branding removed from public entry.\"\nsetParaText(p[9], \"- Landing page Emergent branding removed from public entry.\");\n\n// 2) Section heading: \"New in this update (Branding cleanup)\"\n//    -> \"New in this update (Onboarding flow v1)\"\nsetParaText(p[11], \"New in this update (Onboarding flow v1)\");\n\n// 3) \"- Removed Emergent branding artifacts from frontend public entry:\"\n//    -> \"- Added end-to-end onboarding APIs:\"\nsetParaText(p[12], \"- Added end-to-end onboarding APIs:\");\n\n// 4) \"  - Removed Emergent badge block (`Made with Emergent`).\"\n//    -> \"  - `POST /api/onboarding/start`\"\nsetParaText(p[13], \"  - `POST /api/onboarding/start`\");\n\n// 5) \"  - Removed Emergent external scripts from `index.html`.\"\n//    -> \"  - `POST /api/onboarding/verify-email`\"\nsetParaText(p[14], \"  - `POST /api/onboarding/verify-email`\");\n\n// 6) \"  - Updated page title to `Sitesellr`.\"\n//    -> \"  - `POST /api/onboarding/verify-mobile`\"\nsetParaText(p[15], \"  - `POST /api/onboarding/verify-mobile`\");\n\n// 7) \"  - Updated meta description to `Sitesellr commerce platform`.\"\n//    -> \"  - `GET /api/onboarding/plans`\"\nsetParaText(p[16], \"  - `GET /api/onboarding/plans`\");\n\nawait context.sync();\n\n// 8) Insert seven brand-new paragraphs after the (just-edited) paragraph 16,\n//    each inheriting that paragraph's formatting via insertParagraph.\nconst newLines = [\n  \"  - `POST /api/onboarding/choose-plan`\",\n  \"  - `POST /api/onboarding/confirm-payment` (stub)\",\n  \"  - `POST /api/onboarding/setup-store`\",\n  \"  - `POST /api/onboarding/complete`\",\n  \"- Added frontend onboarding wizard page at `/onboarding` and linked Get Started CTA buttons to it.\",\n  \"- On completion, flow creates User + Merchant + Store + Owner role and signs user in with opaque tokens.\",\n  \"- Onboarding session storage uses in-memory concurrent dictionary (dev-safe, non-persistent).\",\n];\n\nlet anchor = p[16];\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\nawait context.sync();\n\n// Re-load the paragraph collection since the body now has new paragraphs\n// and subsequent original-index-based edits must target the *original*\n// items, which are still valid Word.Paragraph objects (p[] array), so we\n// keep using them directly below.\n\n// 9) \"- Store email/phone/address fields are currently UI-only placeholders\n//     (not persisted in current store model).\"\n//    -> \"- Email/SMS OTP providers are not integrated yet (currently on-screen OTP for dev only).\"\nsetParaText(p[19], \"- Email/SMS OTP providers are not integrated yet (currently on-screen OTP for dev only).\");\n\n// 10) \"- Some advanced form fields are minimal for now (rich order item editor,\n//      full address management UI).\"\n//     -> \"- Payment confirmation is stubbed; no real gateway capture in onboarding yet.\"\nsetParaText(p[20], \"- Payment confirmation is stubbed; no real gateway capture in onboarding yet.\");\n\n// 11) \"- Billing plan enforcement, real payment providers, WebAuthn UI flows, and\n//      production CSP/CORS/rate tuning remain pending.\"\n//     -> \"- Store setup wizard is minimal (no advanced business fields / Cloudflare\n//         provisioning automation yet).\"\nsetParaText(p[21], \"- Store setup wizard is minimal (no advanced business fields / Cloudflare provisioning automation yet).\");\n\nawait context.sync();\n\n// 12) Insert the new \"Billing plan enforcement...\" paragraph (trimmed of\n//     \"real payment providers,\") right after paragraph 21.\np[21].insertParagraph(\n  \"- Billing plan enforcement, WebAuthn UI flows, and production CSP/CORS/rate tuning remain pending.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 13) \"- Full Shopify-level modules (theme marketplace, logistics integrations,\n//      plugin ecosystem, SaaS metering, etc.) remain pending.\"\n//     -> \"- Full Shopify-level modules (themes marketplace, logistics integrations,\n//         plugin ecosystem, SaaS metering, etc.) remain pending.\"\nsetParaText(p[22], \"- Full Shopify-level modules (themes marketplace, logistics integrations, plugin ecosystem, SaaS metering, etc.) remain pending.\");\n\n// 14) \"- Last pushed commit: 9607153\" -> \"- Last pushed commit: cc46cf5\"\nsetParaText(p[25], \"- Last pushed commit: cc46cf5\");\n\n// 15) \"- Current branding cleanup is local and not pushed yet.\"\n//     -> \"- Onboarding flow update is local and pending push.\"\nsetParaText(p[26], \"- Onboarding flow update is local and pending push.\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) \"- Settings General save updates store via backend API.\"\n#    -> \"- Landing page Emergent branding removed from public entry.\"\n$d.Paragraphs.Item(10).Range.Text = \"- Landing page Emergent branding removed from public entry.\"\n\n# 2) Section heading: \"New in this update (Branding cleanup)\"\n#    -> \"New in this update (Onboarding flow v1)\"\n$d.Paragraphs.Item(12).Range.Text = \"New in this update (Onboarding flow v1)\"\n\n# 3) \"- Removed Emergent branding artifacts from frontend public entry:\"\n#    -> \"- Added end-to-end onboarding APIs:\"\n$d.Paragraphs.Item(13).Range.Text = \"- Added end-to-end onboarding APIs:\"\n\n# 4) \"  - Removed Emergent badge block (`Made with Emergent`).\"\n#    -> \"  - `POST /api/onboarding/start`\"\n$d.Paragraphs.Item(14).Range.Text = \"  - ``POST /api/onboarding/start``\"\n\n# 5) \"  - Removed Emergent external scripts from `index.html`.\"\n#    -> \"  - `POST /api/onboarding/verify-email`\"\n$d.Paragraphs.Item(15).Range.Text = \"  - ``POST /api/onboarding/verify-email``\"\n\n# 6) \"  - Updated page title to `Sitesellr`.\"\n#    -> \"  - `POST /api/onboarding/verify-mobile`\"\n$d.Paragraphs.Item(16).Range.Text = \"  - ``POST /api/onboarding/verify-mobile``\"\n\n# 7) \"  - Updated meta description to `Sitesellr commerce platform`.\"\n#    -> \"  - `GET /api/onboarding/plans`\"\n$d.Paragraphs.Item(17).Range.Text = \"  - ``GET /api/onboarding/plans``\"\n\n# 8) Insert seven brand-new paragraphs after paragraph 17, each inheriting\n#    that paragraph's run formatting.\n$newLines = @(\n  \"  - ``POST /api/onboarding/choose-plan``\",\n  \"  - ``POST /api/onboarding/confirm-payment`` (stub)\",\n  \"  - ``POST /api/onboarding/setup-store``\",\n  \"  - ``POST /api/onboarding/complete``\",\n  \"- Added frontend onboarding wizard page at ``/onboarding`` and linked Get Started CTA buttons to it.\",\n  \"- On completion, flow creates User + Merchant + Store + Owner role and signs user in with opaque tokens.\",\n  \"- Onboarding session storage uses in-memory concurrent dictionary (dev-safe, non-persistent).\"\n)\n\n$anchorIndex = 17\nforeach ($line in $newLines) {\n  $anchorPara = $d.Paragraphs.Item($anchorIndex)\n  $anchorPara.Range.InsertParagraphAfter()\n  $anchorIndex = $anchorIndex + 1\n  $d.Paragraphs.Item($anchorIndex).Range.Text = $line\n}\n\n# After inserting 7 paragraphs, every paragraph originally at index >= 18\n# is now shifted down by 7.\n$shift1 = 7\n\n# 9) \"- Store email/phone/address fields are currently UI-only placeholders\n#     (not persisted in current store model).\"\n#    -> \"- Email/SMS OTP providers are not integrated yet (currently on-screen OTP for dev only).\"\n$d.Paragraphs.Item(20 + $shift1).Range.Text = \"- Email/SMS OTP providers are not integrated yet (currently on-screen OTP for dev only).\"\n\n# 10) \"- Some advanced form fields are minimal for now (rich order item editor,\n#      full address management UI).\"\n#     -> \"- Payment confirmation is stubbed; no real gateway capture in onboarding yet.\"\n$d.Paragraphs.Item(21 + $shift1).Range.Text = \"- Payment confirmation is stubbed; no real gateway capture in onboarding yet.\"\n\n# 11) \"- Billing plan enforcement, real payment providers, WebAuthn UI flows, and\n#      production CSP/CORS/rate tuning remain pending.\"\n#     -> \"- Store setup wizard is minimal (no advanced business fields / Cloudflare\n#         provisioning automation yet).\"\n$d.Paragraphs.Item(22 + $shift1).Range.Text = \"- Store setup wizard is minimal (no advanced business fields / Cloudflare provisioning automation yet).\"\n\n# 12) Insert the new \"Billing plan enforcement...\" paragraph (trimmed of\n#     \"real payment providers,\") right after that one.\n$d.Paragraphs.Item(22 + $shift1).Range.InsertParagraphAfter()\n$d.Paragraphs.Item(23 + $shift1).Range.Text = \"- Billing plan enforcement, WebAuthn UI flows, and production CSP/CORS/rate tuning remain pending.\"\n\n$shift2 = 1\n\n# 13) \"- Full Shopify-level modules (theme marketplace, logistics integrations,\n#      plugin ecosystem, SaaS metering, etc.) remain pending.\"\n#     -> \"- Full Shopify-level modules (themes marketplace, logistics integrations,\n#         plugin ecosystem, SaaS metering, etc.) remain pending.\"\n$d.Paragraphs.Item(23 + $shift1 + $shift2).Range.Text = \"- Full Shopify-level modules (themes marketplace, logistics integrations, plugin ecosystem, SaaS metering, etc.) remain pending.\"\n\n# 14) \"- Last pushed commit: 9607153\" -> \"- Last pushed commit: cc46cf5\"\n$d.Paragraphs.Item(26 + $shift1 + $shift2).Range.Text = \"- Last pushed commit: cc46cf5\"\n\n# 15) \"- Current branding cleanup is local and not pushed yet.\"\n#     -> \"- Onboarding flow update is local and pending push.\"\n$d.Paragraphs.Item(27 + $shift1 + $shift2).Range.Text = \"- Onboarding flow update is local and pending push.\"\n"}
